# Auto-generated Excel COM-interop edit script
# Applies value updates to several rows across multiple sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)
# as captured by the source diff for Sheets/Maduin_Profits.xlsx.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()

# Row 58
$ws.Range("H58").Value = 689
$ws.Range("I58").Value = 361.25
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1083.75
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -933.75
$ws.Range("N58").Value = -6300

# Row 86
$ws.Range("H86").Value = 5186.5
$ws.Range("I86").Value = 6132.6665
$ws.Range("J86").Value = 4781
$ws.Range("K86").Value = 6132.6665
$ws.Range("L86").Value = 4781
$ws.Range("M86").Value = -5009.6665
$ws.Range("N86").Value = -7027

# Row 89
$ws.Range("H89").Value = 5186.5
$ws.Range("I89").Value = 6132.6665
$ws.Range("J89").Value = 4781
$ws.Range("K89").Value = 30663.3325
$ws.Range("L89").Value = 23905
$ws.Range("M89").Value = -25047.3325
$ws.Range("N89").Value = -35137

# Row 113
$ws.Range("H113").Value = 21584.857
$ws.Range("I113").Value = 29019.2
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 29019.2
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = -25765.2
$ws.Range("N113").Value = -9507

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 121
$ws.Range("H121").Value = 70000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 70000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 70000
$ws.Range("N121").Value = -73494

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 23
$ws.Range("H23").Value = 75000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 75000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 75000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -75480

# Row 27
$ws.Range("H27").Value = 75000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 75000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 75000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -75384

# Row 31
$ws.Range("H31").Value = 2806.0908
$ws.Range("I31").Value = 2096.889
$ws.Range("J31").Value = 5997.5
$ws.Range("K31").Value = 2096.889
$ws.Range("L31").Value = 5997.5
$ws.Range("M31").Value = -1801.889
$ws.Range("N31").Value = -6587.5

# Row 34
$ws.Range("H34").Value = 2806.0908
$ws.Range("I34").Value = 2096.889
$ws.Range("J34").Value = 5997.5
$ws.Range("K34").Value = 2096.889
$ws.Range("L34").Value = 5997.5
$ws.Range("M34").Value = -1894.889
$ws.Range("N34").Value = -6401.5

# Row 35
$ws.Range("H35").Value = 1634
$ws.Range("I35").Value = 434.85715
$ws.Range("J35").Value = 10028
$ws.Range("K35").Value = 434.85715
$ws.Range("L35").Value = 10028
$ws.Range("M35").Value = -140.85715
$ws.Range("N35").Value = -10616

# Row 38
$ws.Range("H38").Value = 442.25
$ws.Range("I38").Value = 442.25
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 442.25
$ws.Range("M38").Value = -65.25
$ws.Range("N38").ClearContents()

# Row 46
$ws.Range("H46").Value = 442.25
$ws.Range("I46").Value = 442.25
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 442.25
$ws.Range("M46").Value = -231.25
$ws.Range("N46").ClearContents()

# Row 50
$ws.Range("H50").Value = 24997.6
$ws.Range("I50").Value = 25000
$ws.Range("J50").Value = 24997
$ws.Range("K50").Value = 25000
$ws.Range("L50").Value = 24997
$ws.Range("M50").Value = -24375
$ws.Range("N50").Value = -26247

# Row 60
$ws.Range("H60").Value = 16166.667
$ws.Range("I60").Value = 8500

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 44
$ws.Range("H44").Value = 2353.6
$ws.Range("I44").Value = 339.33334
$ws.Range("J44").Value = 5375
$ws.Range("K44").Value = 1018.00002
$ws.Range("L44").Value = 16125
$ws.Range("M44").Value = -620.0000200000001
$ws.Range("N44").Value = -16921

# Row 46
$ws.Range("H46").Value = 2301.4285
$ws.Range("I46").Value = 4100
$ws.Range("J46").Value = 2001.6666
$ws.Range("K46").Value = 12300
$ws.Range("L46").Value = 6004.9998
$ws.Range("M46").Value = -12209
$ws.Range("N46").Value = -6186.9998

# Row 108
$ws.Range("H108").Value = 317
$ws.Range("I108").Value = 317
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 951
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1929

# Row 131
$ws.Range("H131").Value = 982.5
$ws.Range("I131").Value = 900

# Row 140
$ws.Range("H140").Value = 3543.1892
$ws.Range("I140").Value = 736.1818
$ws.Range("J140").Value = 4730.769
$ws.Range("K140").Value = 2208.5454
$ws.Range("L140").Value = 14192.307
$ws.Range("M140").Value = 2971.4546
$ws.Range("N140").Value = -24552.307

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 14
$ws.Range("H14").Value = 166000
$ws.Range("I14").Value = 500450
$ws.Range("J14").Value = 32220
$ws.Range("K14").Value = 500450
$ws.Range("L14").Value = 32220
$ws.Range("M14").Value = -500282
$ws.Range("N14").Value = -32556

# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("N29").ClearContents()

# Row 80
$ws.Range("H80").Value = 7188
$ws.Range("I80").Value = 5631
$ws.Range("J80").Value = 10302
$ws.Range("K80").Value = 5631
$ws.Range("L80").Value = 10302
$ws.Range("M80").Value = -4633
$ws.Range("N80").Value = -12298

# Row 83
$ws.Range("H83").Value = 7188
$ws.Range("I83").Value = 5631
$ws.Range("J83").Value = 10302
$ws.Range("K83").Value = 28155
$ws.Range("L83").Value = 51510
$ws.Range("M83").Value = -23163
$ws.Range("N83").Value = -61494

# Row 102
$ws.Range("H102").Value = 939.26666
$ws.Range("I102").Value = 764.6667
$ws.Range("J102").Value = 1637.6666
$ws.Range("K102").Value = 764.6667
$ws.Range("L102").Value = 1637.6666
$ws.Range("M102").Value = 857.3333
$ws.Range("N102").Value = -4881.6666

# Row 107
$ws.Range("H107").Value = 21750
$ws.Range("I107").Value = 250
$ws.Range("J107").Value = 32500
$ws.Range("K107").Value = 250
$ws.Range("L107").Value = 32500
$ws.Range("M107").Value = 1670
$ws.Range("N107").Value = -36340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 16
$ws.Range("H16").Value = 5450.4
$ws.Range("I16").Value = 4125
$ws.Range("J16").Value = 6334
$ws.Range("K16").Value = 4125
$ws.Range("L16").Value = 6334
$ws.Range("M16").Value = -3955
$ws.Range("N16").Value = -6674

# Row 20
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 1500
$ws.Range("N20").Value = -1952

# Row 23
$ws.Range("H23").Value = 28000000
$ws.Range("I23").Value = 28000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 28000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -27999770

# Row 40
$ws.Range("H40").Value = 1266.1666
$ws.Range("I40").Value = 1248
$ws.Range("J40").Value = 1302.5
$ws.Range("K40").Value = 1248
$ws.Range("L40").Value = 1302.5
$ws.Range("M40").Value = -1112
$ws.Range("N40").Value = -1574.5

# Row 55
$ws.Range("H55").Value = 520.3333
$ws.Range("I55").Value = 480.8889
$ws.Range("J55").Value = 549.9167
$ws.Range("K55").Value = 480.8889
$ws.Range("L55").Value = 549.9167
$ws.Range("M55").Value = -307.8889
$ws.Range("N55").Value = -895.9167

# Row 136
$ws.Range("H136").Value = 4440.154
$ws.Range("I136").Value = 4387.6313
$ws.Range("J136").Value = 4582.7144
$ws.Range("K136").Value = 13162.8939
$ws.Range("L136").Value = 13748.1432
$ws.Range("M136").Value = -10612.8939
$ws.Range("N136").Value = -18848.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 24
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()

# Row 107
$ws.Range("H107").Value = 1104.5385
$ws.Range("I107").Value = 1032.75
$ws.Range("J107").Value = 1219.4
$ws.Range("K107").Value = 3098.25
$ws.Range("L107").Value = 3658.2
$ws.Range("M107").Value = -1178.25
$ws.Range("N107").Value = -7498.200000000001
